$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 315824
$ws.Range("D2").Value = 402612162
$ws.Range("C3").Value = 254
$ws.Range("D3").Value = 302850
$ws.Range("C8").Value = 849
$ws.Range("D8").Value = 1248908
$ws.Range("C10").Value = 115980
$ws.Range("D10").Value = 169951884
$ws.Range("C12").Value = 58525
$ws.Range("D12").Value = 84469011
$ws.Range("C16").Value = 3968
$ws.Range("D16").Value = 5631779
$ws.Range("C19").Value = 69
$ws.Range("D19").Value = 100035
$ws.Range("C20").Value = 6468
$ws.Range("D20").Value = 9024098
$ws.Range("C22").Value = 76354
$ws.Range("D22").Value = 95292750
$ws.Range("C28").Value = 32189
$ws.Range("D28").Value = 47127124
$ws.Range("C30").Value = 11349
$ws.Range("D30").Value = 16324112
$ws.Range("C35").Value = 1777
$ws.Range("D35").Value = 2507333
$ws.Range("C36").Value = 95931
$ws.Range("D36").Value = 120839926
$ws.Range("C42").Value = 898
$ws.Range("D42").Value = 1321685
$ws.Range("C44").Value = 44036
$ws.Range("D44").Value = 64539945
$ws.Range("C46").Value = 9023
$ws.Range("D46").Value = 12949698
$ws.Range("C48").Value = 1393
$ws.Range("D48").Value = 1935139
$ws.Range("C51").Value = 2250
$ws.Range("D51").Value = 3138999
$ws.Range("C52").Value = 68087
$ws.Range("D52").Value = 85456428
$ws.Range("C53").Value = 39
$ws.Range("D53").Value = 43383
$ws.Range("C58").Value = 27889
$ws.Range("D58").Value = 40901626
$ws.Range("C61").Value = 10947
$ws.Range("D61").Value = 15828467
$ws.Range("C63").Value = 1351
$ws.Range("D63").Value = 1888789
$ws.Range("C67").Value = 1438
$ws.Range("D67").Value = 2013599
$ws.Range("C69").Value = 20222
$ws.Range("D69").Value = 26488949
$ws.Range("C73").Value = 7505
$ws.Range("D73").Value = 10987090
$ws.Range("C75").Value = 5044
$ws.Range("D75").Value = 7324206
$ws.Range("C77").Value = 267
$ws.Range("D77").Value = 374173
$ws.Range("C78").Value = 138734
$ws.Range("D78").Value = 173063930
$ws.Range("C79").Value = 68
$ws.Range("D79").Value = 81785
$ws.Range("C84").Value = 62954
$ws.Range("D84").Value = 92274103
$ws.Range("C87").Value = 29328
$ws.Range("D87").Value = 42427589
$ws.Range("C89").Value = 2712
$ws.Range("D89").Value = 3905634
$ws.Range("C90").Value = 2758
$ws.Range("D90").Value = 3897280
$ws.Range("C91").Value = 31993
$ws.Range("D91").Value = 43352309
$ws.Range("C95").Value = 7780
$ws.Range("D95").Value = 11439570
$ws.Range("C97").Value = 7073
$ws.Range("D97").Value = 10253817
$ws.Range("C99").Value = 517
$ws.Range("D99").Value = 734905
$ws.Range("C100").Value = 481
$ws.Range("D100").Value = 693943
$ws.Range("C101").Value = 8747
$ws.Range("D101").Value = 12135796
$ws.Range("C103").Value = 2198
$ws.Range("D103").Value = 3238820
$ws.Range("C105").Value = 2966
$ws.Range("D105").Value = 4331812
$ws.Range("C108").Value = 172
$ws.Range("D108").Value = 244086
$ws.Range("C109").Value = 139231
$ws.Range("D109").Value = 172205766
$ws.Range("C113").Value = 946
$ws.Range("D113").Value = 1389288
$ws.Range("C115").Value = 52206
$ws.Range("D115").Value = 76533682
$ws.Range("C117").Value = 26582
$ws.Range("D117").Value = 38511943
$ws.Range("C121").Value = 2198
$ws.Range("D121").Value = 3087621
$ws.Range("C123").Value = 492992
$ws.Range("D123").Value = 650168871
$ws.Range("C128").Value = 1362
$ws.Range("D128").Value = 2019311
$ws.Range("C130").Value = 204607
$ws.Range("D130").Value = 300785495
$ws.Range("C131").Value = 389
$ws.Range("D131").Value = 580290
$ws.Range("C133").Value = 176750
$ws.Range("D133").Value = 256909711
$ws.Range("C136").Value = 2805
$ws.Range("D136").Value = 3942784
$ws.Range("C138").Value = 6161
$ws.Range("D138").Value = 8705154
$ws.Range("C141").Value = 43753
$ws.Range("D141").Value = 58426527
$ws.Range("C147").Value = 13883
$ws.Range("D147").Value = 20363440
$ws.Range("C148").Value = 3694
$ws.Range("D148").Value = 5328211
$ws.Range("C153").Value = 374
$ws.Range("D153").Value = 527751
$ws.Range("C154").Value = 17223
$ws.Range("D154").Value = 22761148
$ws.Range("C158").Value = 7042
$ws.Range("D158").Value = 10241481
$ws.Range("C160").Value = 4905
$ws.Range("D160").Value = 7060936
$ws.Range("C162").Value = 274
$ws.Range("D162").Value = 378731
$ws.Range("C165").Value = 15315
$ws.Range("D165").Value = 22223846
$ws.Range("C166").Value = 1735
$ws.Range("D166").Value = 2581030
$ws.Range("C171").Value = 86505
$ws.Range("D171").Value = 108219604
$ws.Range("C172").Value = 31
$ws.Range("D172").Value = 36907
$ws.Range("C178").Value = 33528
$ws.Range("D178").Value = 49170947
$ws.Range("C180").Value = 12826
$ws.Range("D180").Value = 18530174
$ws.Range("C182").Value = 1238
$ws.Range("D182").Value = 1732896
$ws.Range("C184").Value = 1607
$ws.Range("D184").Value = 2257433
$ws.Range("C186").Value = 234981
$ws.Range("D186").Value = 292168620
$ws.Range("C187").Value = 135
$ws.Range("D187").Value = 144587
$ws.Range("C194").Value = 85833
$ws.Range("D194").Value = 125827229
$ws.Range("C197").Value = 32596
$ws.Range("D197").Value = 46911532
$ws.Range("C198").Value = 26
$ws.Range("D198").Value = 39000
$ws.Range("C200").Value = 5043
$ws.Range("D200").Value = 7186995
$ws.Range("C203").Value = 4728
$ws.Range("D203").Value = 6542355
$ws.Range("C206").Value = 259829
$ws.Range("D206").Value = 321631088
$ws.Range("C207").Value = 157
$ws.Range("D207").Value = 172473
$ws.Range("C213").Value = 610
$ws.Range("D213").Value = 887878
$ws.Range("C215").Value = 94242
$ws.Range("D215").Value = 137875951
$ws.Range("C218").Value = 50741
$ws.Range("D218").Value = 73331784
$ws.Range("C221").Value = 4619
$ws.Range("D221").Value = 6483916
$ws.Range("C224").Value = 5580
$ws.Range("D224").Value = 7719405
$ws.Range("C227").Value = 104691
$ws.Range("D227").Value = 131043085
$ws.Range("C234").Value = 49038
$ws.Range("D234").Value = 71844270
$ws.Range("C236").Value = 12211
$ws.Range("D236").Value = 17555469
$ws.Range("C238").Value = 1878
$ws.Range("D238").Value = 2691838
$ws.Range("C240").Value = 2429
$ws.Range("D240").Value = 3394138
$ws.Range("C241").Value = 253446
$ws.Range("D241").Value = 320084323
$ws.Range("C249").Value = 94772
$ws.Range("D249").Value = 138875512
$ws.Range("C252").Value = 63995
$ws.Range("D252").Value = 92736143
$ws.Range("C254").Value = 2384
$ws.Range("D254").Value = 3364750
$ws.Range("C257").Value = 4488
$ws.Range("D257").Value = 6300040
